$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "Type"
$ws.Range("B1").Value = "Weight"

# --- Column A (Type) for rows 2-7, in top-to-bottom order so new shared
#     strings are appended in the same order Excel would record them ---
$ws.Range("A2").Value = "SN1_[FA-H]-"
$ws.Range("A3").Value = "SN2_[FA-H]-"
$ws.Range("A4").Value = "[LPL(SN1)-H]-"
$ws.Range("A5").Value = "[LPL(SN2)-H]-"
$ws.Range("A6").Value = "[LPL(SN1)-H2O-H]-"
$ws.Range("A7").Value = "[LPL(SN2)-H2O-H]-"

# --- New Group header, written after column A so it lands last in the
#     shared-string table (matches the target uniqueCount ordering) ---
$ws.Range("C1").Value = "Group"

# --- Column B (Weight) values (unchanged from before, but set explicitly) ---
$ws.Range("B2").Value = 40
$ws.Range("B3").Value = 40
$ws.Range("B4").Value = 5
$ws.Range("B5").Value = 5
$ws.Range("B6").Value = 5
$ws.Range("B7").Value = 5

# --- Column C (Group) numeric values, replacing former C:F data ---
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 2
$ws.Range("C6").Value = 2
$ws.Range("C7").Value = 2

# --- Remove the old D:F columns' content entirely (now unused) ---
$ws.Range("D1:F7").Clear()

# --- Selection / view state to match the saved workbook ---
$ws.Range("D1:F1048576").Select()
